$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.445647641019636, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 5.507293877332936)
    3 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 17.08608867836142)
    4 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    5 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    6 = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.038307959104277)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
